$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 is the "quartz scheduler" project row. The commit adds the
# dynamic/static/grep analysis (plus their notes) for quartz, and fixes
# the download link.

# D21: link text - point at the downloads page instead of the bare domain.
$ws.Range("D21").Value = "http://quartz-scheduler.org/downloads"

# E21: Dynamic Analysis
$ws.Range("E21").Value = "Callback on VMStart.
#################################################
SecurityManager Changed:
QuartzServer.java, main, 178
#################################################
Agent OnUnload, agent exits."

# F21: Static Analysis
$ws.Range("F21").Value = "At QuartzServer.java:[line 178]
In method org.quartz.impl.QuartzServer.main(String[])
Value Not null: java.rmi.RMISecurityManager
Value new 
Value new[187](3) 37 
Value Variable is set at:"

# G21: Grep Results
$ws.Range("G21").Value = "./src/org/quartz/impl/QuartzServer.java
177:        if (System.getSecurityManager() == null) {
178:            System.setSecurityManager(new java.rmi.RMISecurityManager());"

# H21: General Notes
$ws.Range("H21").Value = "The program sets an RMISecurityManager if no SecurityManager is set when running the main program in QuartzServer.java.  Other than that no interaction with the SecurityManager occurs (so no nulling or weakening happening)."

# I21: Dynamic Analysis Notes
$ws.Range("I21").Value = "Sets SecurityManager at start of program as expected from Grep results"

# J21: Static Analysis Notes
$ws.Range("J21").Value = "At QuartzServer.java:[line 178]
In method org.quartz.impl.QuartzServer.main(String[])
Value Not null: java.rmi.RMISecurityManager
Value new 
Value new[187](3) 37 
Value Variable is set at:
Sets the SecurityManager at the beginning of the QuartzServer startup if the program is not started with a SecurityManager already running."

# K21: Grep Notes
$ws.Range("K21").Value = "./src/org/quartz/impl/QuartzServer.java
177:        if (System.getSecurityManager() == null) {
178:            System.setSecurityManager(new java.rmi.RMISecurityManager());
First two lines of the main function.  If the program is ran without a SecurityManager, set the SecurityManager to the RMISecurityManager.
"

# Move / record the active selection as it ended up after this edit (C20).
[void]$ws.Range("C20").Select()
